# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the Leve profit tables across all eight job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner price update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 40426
$ws.Range("J87").Value = 40426
$ws.Range("L87").Value = 40426
$ws.Range("N87").Value = -42922

$ws.Range("H90").Value = 40426
$ws.Range("J90").Value = 40426
$ws.Range("L90").Value = 121278
$ws.Range("N90").Value = -133758

$ws.Range("H111").Value = 3461.2856
$ws.Range("I111").Value = 7143
$ws.Range("J111").Value = 700
$ws.Range("K111").Value = 21429
$ws.Range("L111").Value = 2100
$ws.Range("M111").Value = -18362
$ws.Range("N111").Value = -8234

$ws.Range("H132").Value = 6543555
$ws.Range("I132").Value = 7412651.5
$ws.Range("K132").Value = 22237954.5
$ws.Range("M132").Value = -22235424.5

$ws.Range("H137").Value = 2246.2
$ws.Range("I137").Value = 3222.25
$ws.Range("K137").Value = 9666.75
$ws.Range("M137").Value = -7116.75

$ws.Range("H138").Value = 3456.5764
$ws.Range("I138").Value = 1895.3684
$ws.Range("J138").Value = 3906.0151
$ws.Range("K138").Value = 5686.1052
$ws.Range("L138").Value = 11718.0453
$ws.Range("M138").Value = -546.1052
$ws.Range("N138").Value = -21998.0453

$ws.Range("H141").Value = 1293.9412
$ws.Range("I141").Value = 1206.3846
$ws.Range("K141").Value = 3619.1538
$ws.Range("M141").Value = 1560.8462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20448.678
$ws.Range("I32").Value = 18560.154
$ws.Range("J32").Value = 44999.5
$ws.Range("K32").Value = 18560.154
$ws.Range("L32").Value = 44999.5
$ws.Range("M32").Value = -18273.154
$ws.Range("N32").Value = -45573.5

$ws.Range("H61").Value = 34484184
$ws.Range("I61").Value = 41667908
$ws.Range("J61").Value = 2308.4
$ws.Range("K61").Value = 41667908
$ws.Range("L61").Value = 2308.4
$ws.Range("M61").Value = -41667696
$ws.Range("N61").Value = -2732.4

$ws.Range("H74").Value = 1582.2222
$ws.Range("I74").Value = 621.95
$ws.Range("K74").Value = 621.95
$ws.Range("M74").Value = 252.05

$ws.Range("H77").Value = 1582.2222
$ws.Range("I77").Value = 621.95
$ws.Range("K77").Value = 3109.75
$ws.Range("M77").Value = 1258.25

$ws.Range("H110").Value = 676.5
$ws.Range("J110").Value = 1100
$ws.Range("L110").Value = 1100
$ws.Range("N110").Value = -5190

$ws.Range("H122").Value = 4388.1113
$ws.Range("I122").Value = 4272.2
$ws.Range("J122").Value = 4533
$ws.Range("K122").Value = 12816.6
$ws.Range("L122").Value = 13599
$ws.Range("M122").Value = -10366.6
$ws.Range("N122").Value = -18499

$ws.Range("H132").Value = 2682.718
$ws.Range("I132").Value = 1862.2174
$ws.Range("J132").Value = 3862.1875
$ws.Range("K132").Value = 5586.6522
$ws.Range("L132").Value = 11586.5625
$ws.Range("M132").Value = -3056.6522
$ws.Range("N132").Value = -16646.5625

$ws.Range("H136").Value = 34484184
$ws.Range("I136").Value = 41667908
$ws.Range("J136").Value = 2308.4
$ws.Range("K136").Value = 125003724
$ws.Range("L136").Value = 6925.200000000001
$ws.Range("M136").Value = -125001174
$ws.Range("N136").Value = -12025.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3487.9546
$ws.Range("I134").Value = 939.88464
$ws.Range("K134").Value = 2819.65392
$ws.Range("M134").Value = -284.6539199999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1606.9524
$ws.Range("I31").Value = 1381.5438
$ws.Range("J31").Value = 3748.3333
$ws.Range("K31").Value = 1381.5438
$ws.Range("L31").Value = 3748.3333
$ws.Range("M31").Value = -1086.5438
$ws.Range("N31").Value = -4338.3333

$ws.Range("H34").Value = 1606.9524
$ws.Range("I34").Value = 1381.5438
$ws.Range("J34").Value = 3748.3333
$ws.Range("K34").Value = 1381.5438
$ws.Range("L34").Value = 3748.3333
$ws.Range("M34").Value = -1179.5438
$ws.Range("N34").Value = -4152.3333

$ws.Range("H50").Value = 25000
$ws.Range("J50").Value = 25000
$ws.Range("L50").Value = 25000
$ws.Range("N50").Value = -26250

$ws.Range("H51").Value = 25000
$ws.Range("J51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -26472

$ws.Range("H58").Value = 3553.848
$ws.Range("I58").Value = 1012.94446
$ws.Range("K58").Value = 1012.94446
$ws.Range("M58").Value = -809.94446

$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H61").Value = 25000
$ws.Range("J61").Value = 25000
$ws.Range("L61").Value = 25000
$ws.Range("N61").Value = -25696

$ws.Range("H134").Value = 12196719
$ws.Range("I134").Value = 1521.25
$ws.Range("K134").Value = 4563.75
$ws.Range("M134").Value = -2028.75

$ws.Range("H136").Value = 3553.848
$ws.Range("I136").Value = 1012.94446
$ws.Range("K136").Value = 3038.83338
$ws.Range("M136").Value = -488.83338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 59.666668
$ws.Range("I40").Value = 65.875
$ws.Range("J40").Value = 10
$ws.Range("K40").Value = 263.5
$ws.Range("L40").Value = 40
$ws.Range("M40").Value = -194.5
$ws.Range("N40").Value = -178

$ws.Range("H140").Value = 29923.162
$ws.Range("J140").Value = 2889.5217
$ws.Range("L140").Value = 8668.5651
$ws.Range("N140").Value = -19028.5651

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 598.1739
$ws.Range("I97").Value = 576.9231
$ws.Range("J97").Value = 625.8
$ws.Range("K97").Value = 576.9231
$ws.Range("L97").Value = 625.8
$ws.Range("M97").Value = -80.92309999999998
$ws.Range("N97").Value = -1617.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1337.1428
$ws.Range("I46").Value = 453.33334
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 453.33334
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -265.33334
$ws.Range("N46").Value = -2376

$ws.Range("H82").Value = 2045.762
$ws.Range("I82").Value = 2073.5
$ws.Range("J82").Value = 2008.7778
$ws.Range("K82").Value = 2073.5
$ws.Range("L82").Value = 2008.7778
$ws.Range("M82").Value = -1712.5
$ws.Range("N82").Value = -2730.7778

$ws.Range("H85").Value = 2045.762
$ws.Range("I85").Value = 2073.5
$ws.Range("J85").Value = 2008.7778
$ws.Range("K85").Value = 2073.5
$ws.Range("L85").Value = 2008.7778
$ws.Range("M85").Value = -825.5
$ws.Range("N85").Value = -4504.7778

$ws.Range("H132").Value = 2779.158
$ws.Range("I132").Value = 2080.2
$ws.Range("J132").Value = 3555.7778
$ws.Range("K132").Value = 6240.599999999999
$ws.Range("L132").Value = 10667.3334
$ws.Range("M132").Value = -3710.599999999999
$ws.Range("N132").Value = -15727.3334

$ws.Range("H136").Value = 1906.75
$ws.Range("I136").Value = 1725.3889
$ws.Range("J136").Value = 2450.8333
$ws.Range("K136").Value = 5176.1667
$ws.Range("L136").Value = 7352.499899999999
$ws.Range("M136").Value = -2626.1667
$ws.Range("N136").Value = -12452.4999

$ws.Range("H139").Value = 35580
$ws.Range("J139").Value = 35580
$ws.Range("L139").Value = 35580
$ws.Range("N139").Value = -45860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3744
$ws.Range("I132").Value = 3840.2666
$ws.Range("K132").Value = 11520.7998
$ws.Range("M132").Value = -8990.799800000001
